$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column 54 (BB) width from 11.7109375 to 12.7109375 (closest achievable via COM pixel-quantized width)
$ws.Columns.Item(54).ColumnWidth = 11.83

# Update cell values per diff
$ws.Range("AN1").Value = 0.67492322814036076
$ws.Range("N2").Value = 0.93016850520602001
$ws.Range("A3").Value = 0.86983215623352339
$ws.Range("N3").Value = 0.90242982680627937
$ws.Range("S3").Value = 0.81360932073736869
$ws.Range("BA3").Value = 0.67530354170618356
$ws.Range("B4").Value = 0.87516004882373788
$ws.Range("AT4").Value = 0.97947252836007315
$ws.Range("E6").Value = 0.83416521379654984
$ws.Range("S6").Value = 0.92693425463300638
$ws.Range("E7").Value = 0.69718521129719635
$ws.Range("F7").Value = 0.81944690911522888
$ws.Range("I7").Value = 0.97831146113056855
$ws.Range("BP7").Value = 0.77534025101199822
$ws.Range("AE8").Value = 0.54829804251441483
$ws.Range("AP8").Value = 0.89146884196465703
$ws.Range("BD8").Value = 0.75582725638848725
$ws.Range("J9").Value = 0.90503280693923616
$ws.Range("AX9").Value = 0.93719110205514955
$ws.Range("B10").Value = 0.71538302171705959
$ws.Range("A11").Value = 0.86972655137541754
$ws.Range("H11").Value = 0.95431048619355963
$ws.Range("AC11").Value = 0.97436995381872871
$ws.Range("BC11").Value = 0.8185471720861357
$ws.Range("AZ12").Value = 0.52399023924596166
$ws.Range("Q13").Value = 0.97380478442949203
$ws.Range("AT13").Value = 0.7993264005018389
$ws.Range("AP14").Value = 0.74326060495803192
$ws.Range("O16").Value = 0.88835328853097151
$ws.Range("BK16").Value = 0.80815189014150013
$ws.Range("D17").Value = 0.67428410751203705
$ws.Range("J17").Value = 0.68565628067660378
$ws.Range("T18").Value = 0.8395158836021408
$ws.Range("BN18").Value = 0.77128440207991567
$ws.Range("Q19").Value = 0.73265277829846942
$ws.Range("U19").Value = 0.69752382143590264
$ws.Range("AK20").Value = 0.71749747966200161
$ws.Range("X21").Value = 0.96463023323262043
$ws.Range("J22").Value = 0.83744598333469011
$ws.Range("T22").Value = 0.94719141374006954
$ws.Range("W22").Value = 0.81183470206900776
$ws.Range("BF23").Value = 0.82413060369091018
$ws.Range("A24").Value = 0.80121732854210426
$ws.Range("O24").Value = 0.97897034956611706
$ws.Range("BF24").Value = 0.83112516834316186
$ws.Range("BH25").Value = 0.70893506597599187
$ws.Range("H26").Value = 0.99079816501995799
$ws.Range("AF27").Value = 0.77190422161971672
$ws.Range("AK27").Value = 0.8555702130978966
$ws.Range("G28").Value = 0.7619059555178892
$ws.Range("O28").Value = 0.78810739392061946
$ws.Range("AD28").Value = 0.87366438954689118
$ws.Range("J29").Value = 0.89497814921078689
$ws.Range("AE29").Value = 0.69956251740254571
$ws.Range("BO29").Value = 0.98241224543835293
$ws.Range("B30").Value = 0.94301682275884291
$ws.Range("S30").Value = 0.85333736957547679
$ws.Range("AA30").Value = 0.86516027514691785
$ws.Range("BO30").Value = 0.9117325021259739
$ws.Range("U32").Value = 0.63444883560941645
$ws.Range("AE32").Value = 0.60687077316847338
$ws.Range("AM32").Value = 0.9086831713088892
$ws.Range("BK32").Value = 0.85654480180969705
$ws.Range("N33").Value = 0.5442057814876049
$ws.Range("AI33").Value = 0.6709593325953418
$ws.Range("BI33").Value = 0.92533403334355979
$ws.Range("I34").Value = 0.54003652274311476
$ws.Range("N34").Value = 0.94872935182820917
$ws.Range("AB34").Value = 0.78648104469434355
$ws.Range("AZ34").Value = 0.93120194978999082
$ws.Range("Y35").Value = 0.84837989870114172
$ws.Range("AU35").Value = 0.85279531926072005
$ws.Range("AZ35").Value = 0.66961535073228506
$ws.Range("T36").Value = 0.87701199155934262
$ws.Range("P38").Value = 0.98855163182275863
$ws.Range("AS38").Value = 0.71198915088135883
$ws.Range("Y40").Value = 0.99024905698219556
$ws.Range("BB40").Value = 0.94913601228700628
$ws.Range("BC40").Value = 0.69582812680473505
$ws.Range("AP41").Value = 0.97860629116354469
$ws.Range("BE41").Value = 0.72862474366717411
$ws.Range("D43").Value = 0.93559708018305443
$ws.Range("F43").Value = 0.76209844076561728
$ws.Range("AJ43").Value = 0.9480963733393899
$ws.Range("AO43").Value = 0.67056289994852925
$ws.Range("AJ44").Value = 0.7047393025086679
$ws.Range("BF44").Value = 0.96791281499422799
$ws.Range("U45").Value = 0.63631749951362604
$ws.Range("AM45").Value = 0.86556123583567168
$ws.Range("AU45").Value = 0.69288468085332799
$ws.Range("C46").Value = 0.87713591886966102
$ws.Range("AA46").Value = 0.80150988493248454
$ws.Range("BP47").Value = 0.78878240354619567
$ws.Range("AX48").Value = 0.83519713185661071
$ws.Range("BF48").Value = 0.88500615351785483
$ws.Range("A49").Value = 0.79101317257061221
$ws.Range("BK49").Value = 0.96621913316907904
$ws.Range("F50").Value = 0.88380130725799144
$ws.Range("P50").Value = 0.98496688355997386
$ws.Range("AF50").Value = 0.71993598558798644
$ws.Range("AI51").Value = 0.95369623002146375
$ws.Range("AT51").Value = 0.83256890732943223
$ws.Range("AW51").Value = 0.93842866893108723
$ws.Range("BA52").Value = 0.77545227335396238
$ws.Range("BE53").Value = 0.82654632838082565
$ws.Range("I54").Value = 0.99371951599938457
$ws.Range("W56").Value = 0.69640806963277768
$ws.Range("BF56").Value = 0.95501322308165604
$ws.Range("BJ56").Value = 0.87616964602070602
$ws.Range("X57").Value = 0.99822803363218149
$ws.Range("AA57").Value = 0.67203756632520228
$ws.Range("L58").Value = 0.89528319199823458
$ws.Range("Y59").Value = 0.99375630229887213
$ws.Range("Z59").Value = 0.8204119112387982
$ws.Range("BI59").Value = 0.83174211459417657
$ws.Range("K60").Value = 0.66098487110262516
$ws.Range("AV60").Value = 0.78209214759881329
$ws.Range("BG60").Value = 0.58850748418615639
$ws.Range("G61").Value = 0.91707995512488916
$ws.Range("N61").Value = 0.86313756051540214
$ws.Range("Y62").Value = 0.78863581123243032
$ws.Range("X63").Value = 0.8187654159481248
$ws.Range("AT63").Value = 0.81907690243718001
$ws.Range("M64").Value = 0.93000247104385014
$ws.Range("AK64").Value = 0.91034113916764736
$ws.Range("N65").Value = 0.85298730329269379
$ws.Range("T65").Value = 0.73396524504089156
$ws.Range("AY65").Value = 0.80637157790085601
$ws.Range("AN66").Value = 0.95322077067574762
$ws.Range("F67").Value = 0.83527658132891192
$ws.Range("V67").Value = 0.97125648808967613
$ws.Range("AP67").Value = 0.93440077814147493
$ws.Range("AV67").Value = 0.83796058983337973
$ws.Range("BN67").Value = 0.65479360050004076
$ws.Range("Y68").Value = 0.73966117953056587
$ws.Range("BH68").Value = 0.92389864086526785
